# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note text (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newLine = "`n"
$text = "Conversión del día 💰" + $newLine +
        "✅ Dólar paralelo: 68" + $newLine +
        $newLine +
        "Binance" + $newLine +
        "✅ 1000 Bs = 3.37 = 13045.7 pesos" + $newLine +
        "✅ 13045.7 pesos = 3.35 = 965.85 Bs" + $newLine +
        $newLine +
        "Promedio competencia" + $newLine +
        "✅ Tasa pesos: 20" + $newLine +
        "✅ Tasa Bs: 20" + $newLine +
        "✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $text

# --- tasas: update the rate cells N10, O10, N12, O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 296.5
$wsTasas.Range("O10").Value = 3868.05
$wsTasas.Range("N12").Value = 3890
$wsTasas.Range("O12").Value = 288
